# "aggiornamento a 9/09 compreso" - append daily rows through 2021-09-09
# (Excel serial dates 44441..44448) to the bottom of the data table,
# continuing the existing A:D layout (date | nuovi pos. | somma mobile 7gg. |
# somma mobile 7gg. per 100mila abitanti).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count  # 366 -> new rows start at 367

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile per 100k)
$data = @(
    @(44441, 1, 14, 139.0268123138034),
    @(44442, 0, 12, 119.1658391261172),
    @(44443, 4, 16, 158.8877855014896),
    @(44444, 0, 9, 89.37437934458789),
    @(44445, 2, 9, 89.37437934458789),
    @(44446, 0, 8, 79.44389275074478),
    @(44447, 0, 7, 69.51340615690168),
    @(44448, 0, 6, 59.5829195630586)
)

$firstNewRow = $lastRow + 1
$lastNewRow = $lastRow + $data.Count

# Copy the formatting (date style with s="2" on column A, etc.) from the
# last existing row down onto the new rows before filling in values.
$ws.Range("A$lastRow`:D$lastRow").Copy()
$ws.Range("A$firstNewRow`:D$lastNewRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$r = $firstNewRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
